# ----------------------------------------------------------------------------
# Scheduled runner update: refresh market-board derived columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
#  LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ) for the Leve
# profitability sheets, pulling the latest market-board snapshot values.
# ----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 18
$ws.Cells.Item(18, 8).Value = 892.44446
$ws.Cells.Item(18, 9).Value = 892.44446
$ws.Cells.Item(18, 11).Value = 892.44446
$ws.Cells.Item(18, 13).Value = -608.44446

# Row 33
$ws.Cells.Item(33, 8).Value = 495.5
$ws.Cells.Item(33, 9).Value = 495.5
$ws.Cells.Item(33, 11).Value = 495.5
$ws.Cells.Item(33, 13).Value = -266.5

# Row 93
$ws.Cells.Item(93, 8).Value = 818363700
$ws.Cells.Item(93, 10).Value = 1000000
$ws.Cells.Item(93, 12).Value = 1000000
$ws.Cells.Item(93, 14).Value = -1004992

# Row 95
$ws.Cells.Item(95, 8).Value = 36000
$ws.Cells.Item(95, 10).Value = 36000
$ws.Cells.Item(95, 12).Value = 36000
$ws.Cells.Item(95, 14).Value = -41492

# Row 100
$ws.Cells.Item(100, 8).Value = 3862.75
$ws.Cells.Item(100, 9).Value = 3181.4
$ws.Cells.Item(100, 11).Value = 3181.4
$ws.Cells.Item(100, 13).Value = -2640.4

# Row 138
$ws.Cells.Item(138, 8).Value = 4071.7104
$ws.Cells.Item(138, 10).Value = 4155.9
$ws.Cells.Item(138, 12).Value = 12467.7
$ws.Cells.Item(138, 14).Value = -22747.7

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 3964.3125
$ws.Cells.Item(32, 9).Value = 3339.8
$ws.Cells.Item(32, 10).Value = 13332
$ws.Cells.Item(32, 11).Value = 3339.8
$ws.Cells.Item(32, 12).Value = 13332
$ws.Cells.Item(32, 13).Value = -3052.8
$ws.Cells.Item(32, 14).Value = -13906

# Row 76
$ws.Cells.Item(76, 8).Value = 32486.25
$ws.Cells.Item(76, 10).Value = 32486.25
$ws.Cells.Item(76, 12).Value = 32486.25
$ws.Cells.Item(76, 14).Value = -33162.25

# Row 79
$ws.Cells.Item(79, 8).Value = 32486.25
$ws.Cells.Item(79, 10).Value = 32486.25
$ws.Cells.Item(79, 12).Value = 32486.25
$ws.Cells.Item(79, 14).Value = -34826.25

# Row 95
$ws.Cells.Item(95, 8).Value = 0
$ws.Cells.Item(95, 10).Value = 0
$ws.Cells.Item(95, 12).Value = 0
$ws.Cells.Item(95, 14).ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 64
$ws.Cells.Item(64, 8).Value = 891.3333
$ws.Cells.Item(64, 9).Value = 893
$ws.Cells.Item(64, 10).Value = 888
$ws.Cells.Item(64, 11).Value = 893
$ws.Cells.Item(64, 12).Value = 888
$ws.Cells.Item(64, 13).Value = -668
$ws.Cells.Item(64, 14).Value = -1338

# Row 67
$ws.Cells.Item(67, 8).Value = 891.3333
$ws.Cells.Item(67, 9).Value = 893
$ws.Cells.Item(67, 10).Value = 888
$ws.Cells.Item(67, 11).Value = 893
$ws.Cells.Item(67, 12).Value = 888
$ws.Cells.Item(67, 13).Value = -113
$ws.Cells.Item(67, 14).Value = -2448

# Row 75
$ws.Cells.Item(75, 8).Value = 6107
$ws.Cells.Item(75, 9).Value = 6107
$ws.Cells.Item(75, 11).Value = 6107
$ws.Cells.Item(75, 13).Value = -5171

# Row 78
$ws.Cells.Item(78, 8).Value = 6107
$ws.Cells.Item(78, 9).Value = 6107
$ws.Cells.Item(78, 11).Value = 18321
$ws.Cells.Item(78, 13).Value = -13641

# Row 86
$ws.Cells.Item(86, 8).Value = 2145.25
$ws.Cells.Item(86, 9).Value = 1473.2
$ws.Cells.Item(86, 10).Value = 4161.4
$ws.Cells.Item(86, 11).Value = 1473.2
$ws.Cells.Item(86, 12).Value = 4161.4
$ws.Cells.Item(86, 13).Value = -350.2
$ws.Cells.Item(86, 14).Value = -6407.4

# Row 89
$ws.Cells.Item(89, 8).Value = 2145.25
$ws.Cells.Item(89, 9).Value = 1473.2
$ws.Cells.Item(89, 10).Value = 4161.4
$ws.Cells.Item(89, 11).Value = 7366
$ws.Cells.Item(89, 12).Value = 20807
$ws.Cells.Item(89, 13).Value = -1750
$ws.Cells.Item(89, 14).Value = -32039

# Row 134
$ws.Cells.Item(134, 8).Value = 2419.4443
$ws.Cells.Item(134, 10).Value = 3333.3333
$ws.Cells.Item(134, 12).Value = 9999.999899999999
$ws.Cells.Item(134, 14).Value = -15069.9999

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Cells.Item(7, 8).Value = 199.05263
$ws.Cells.Item(7, 9).Value = 234.71428
$ws.Cells.Item(7, 11).Value = 234.71428
$ws.Cells.Item(7, 13).Value = -121.71428

# Row 132
$ws.Cells.Item(132, 8).Value = 1181.4
$ws.Cells.Item(132, 9).Value = 976.75
$ws.Cells.Item(132, 11).Value = 2930.25
$ws.Cells.Item(132, 13).Value = -400.25

$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Cells.Item(68, 8).Value = 2926.5715
$ws.Cells.Item(68, 9).Value = 2996.75
$ws.Cells.Item(68, 10).Value = 2898.5
$ws.Cells.Item(68, 11).Value = 8990.25
$ws.Cells.Item(68, 12).Value = 8695.5
$ws.Cells.Item(68, 13).Value = -8179.25
$ws.Cells.Item(68, 14).Value = -10317.5

# Row 71
$ws.Cells.Item(71, 8).Value = 2926.5715
$ws.Cells.Item(71, 9).Value = 2996.75
$ws.Cells.Item(71, 10).Value = 2898.5
$ws.Cells.Item(71, 11).Value = 26970.75
$ws.Cells.Item(71, 12).Value = 26086.5
$ws.Cells.Item(71, 13).Value = -22914.75
$ws.Cells.Item(71, 14).Value = -34198.5

# Row 107
$ws.Cells.Item(107, 8).Value = 239.66667
$ws.Cells.Item(107, 9).Value = 239.66667
$ws.Cells.Item(107, 11).Value = 719.00001
$ws.Cells.Item(107, 13).Value = 1200.99999

# Row 131
$ws.Cells.Item(131, 8).Value = 1000.8889
$ws.Cells.Item(131, 10).Value = 3000
$ws.Cells.Item(131, 12).Value = 9000
$ws.Cells.Item(131, 14).Value = -19080

# Row 133
$ws.Cells.Item(133, 8).Value = 0
$ws.Cells.Item(133, 9).Value = 0
$ws.Cells.Item(133, 11).Value = 0
$ws.Cells.Item(133, 13).ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 105
$ws.Cells.Item(105, 8).Value = 12249
$ws.Cells.Item(105, 10).Value = 12249
$ws.Cells.Item(105, 12).Value = 12249
$ws.Cells.Item(105, 14).Value = -19237

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Cells.Item(22, 8).Value = 923.14636
$ws.Cells.Item(22, 10).Value = 1206.125
$ws.Cells.Item(22, 12).Value = 1206.125
$ws.Cells.Item(22, 14).Value = -1796.125

# Row 27
$ws.Cells.Item(27, 8).Value = 923.14636
$ws.Cells.Item(27, 10).Value = 1206.125
$ws.Cells.Item(27, 12).Value = 1206.125
$ws.Cells.Item(27, 14).Value = -1420.125

$ws = $wb.Worksheets.Item("WVR")
# Row 69
$ws.Cells.Item(69, 8).Value = 32004.4
$ws.Cells.Item(69, 10).Value = 32004.4
$ws.Cells.Item(69, 12).Value = 32004.4
$ws.Cells.Item(69, 14).Value = -33502.4

# Row 72
$ws.Cells.Item(72, 8).Value = 32004.4
$ws.Cells.Item(72, 10).Value = 32004.4
$ws.Cells.Item(72, 12).Value = 96013.20000000001
$ws.Cells.Item(72, 14).Value = -103501.2

# Row 80
$ws.Cells.Item(80, 8).Value = 29099.8
$ws.Cells.Item(80, 10).Value = 29099.8
$ws.Cells.Item(80, 12).Value = 29099.8
$ws.Cells.Item(80, 14).Value = -31095.8

# Row 83
$ws.Cells.Item(83, 8).Value = 29099.8
$ws.Cells.Item(83, 10).Value = 29099.8
$ws.Cells.Item(83, 12).Value = 87299.39999999999
$ws.Cells.Item(83, 14).Value = -97283.39999999999

# Row 126
$ws.Cells.Item(126, 8).Value = 4108.909
$ws.Cells.Item(126, 9).Value = 1911.625
$ws.Cells.Item(126, 10).Value = 9968.333000000001
$ws.Cells.Item(126, 11).Value = 5734.875
$ws.Cells.Item(126, 12).Value = 29904.999
$ws.Cells.Item(126, 13).Value = -3264.875
$ws.Cells.Item(126, 14).Value = -34844.999

# Row 132
$ws.Cells.Item(132, 8).Value = 4298.8335
$ws.Cells.Item(132, 9).Value = 3948.5
$ws.Cells.Item(132, 10).Value = 4999.5
$ws.Cells.Item(132, 11).Value = 11845.5
$ws.Cells.Item(132, 12).Value = 14998.5
$ws.Cells.Item(132, 13).Value = -9315.5
$ws.Cells.Item(132, 14).Value = -20058.5
